$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")

$ws.Range("A2").Value = "Última actualización: 06:49:33"
$ws.Range("A3").Value = "Total filas: 73"

$rows = @(
    ,@(47, "06:49:33", "06:59", "14_ABASTO", 10, "LP1912")
    ,@(49, "06:49:33", "07:01", "16_SANTA ANA", 12, "LP1912")
    ,@(50, "06:49:33", "07:04", "23_HERNANDEZ", 15, "LP1912")
    ,@(51, "06:49:33", "07:05", "15_ABASTO", 16, "LP1912")
    ,@(52, "05:52:07", "07:05", "23_HERNANDEZ", 73, "LP1912")
    ,@(53, "06:49:33", "07:07", "225_GOMEZ", 18, "LP1912")
    ,@(54, "06:49:33", "07:11", "215A_EL PATO", 22, "LP1912")
    ,@(55, "05:52:07", "07:12", "215A_EL PATO", 80, "LP1912")
    ,@(56, "06:49:33", "07:15", "11_ETCHEVERRY", 26, "LP1912")
    ,@(57, "06:49:33", "07:16", "16_SANTA ANA", 27, "LP1912")
    ,@(58, "05:52:07", "07:16", "11_ETCHEVERRY", 84, "LP1912")
    ,@(59, "06:49:33", "07:21", "26_HERNANDEZ", 32, "LP1912")
    ,@(60, "06:49:33", "07:23", "10_OLMOS", 34, "LP1912")
    ,@(61, "06:49:33", "07:31", "11_ETCHEVERRY", 42, "LP1912")
    ,@(62, "05:52:07", "07:32", "16_SANTA ANA", 100, "LP1912")
    ,@(63, "06:49:33", "07:32", "84_COLONIA URQUIZA-ESC 49", 43, "LP1912")
    ,@(64, "05:52:07", "07:32", "11_ETCHEVERRY", 100, "LP1912")
    ,@(65, "06:49:33", "07:36", "27_EL RETIRO", 47, "LP1912")
    ,@(66, "06:21:22", "07:37", "27_EL RETIRO", 76, "LP1912")
    ,@(67, "06:49:33", "07:39", "10_OLMOS", 50, "LP1912")
    ,@(68, "06:49:33", "07:47", "14_ABASTO", 58, "LP1912")
    ,@(69, "05:52:07", "07:48", "14_ABASTO", 116, "LP1912")
    ,@(70, "06:49:33", "07:51", "215D_EL PATO", 62, "LP1912")
    ,@(71, "06:21:22", "08:01", "23_HERNANDEZ", 100, "LP1912")
    ,@(72, "06:49:33", "08:06", "23_HERNANDEZ", 77, "LP1912")
    ,@(73, "06:49:33", "08:12", "15_ABASTO", 83, "LP1912")
    ,@(74, "06:49:33", "08:21", "26_HERNANDEZ", 92, "LP1912")
    ,@(75, "06:49:33", "08:22", "16_P MOR-SANTA ANA", 93, "LP1912")
    ,@(76, "06:49:33", "08:23", "215B_EL PATO", 94, "LP1912")
    ,@(77, "06:49:33", "08:27", "84_COLONIA URQUIZA-ESC 49", 98, "LP1912")
    ,@(78, "06:49:33", "08:42", "81_EL PELIGRO", 113, "LP1912")
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")

$ws.Range("A2").Value = "Última actualización: 06:49:33"
$ws.Range("A3").Value = "Total filas: 16"

$rows = @(
    ,@(18, "06:49:33", "07:11", "215A_EL PATO", 22, "LP1912")
    ,@(20, "06:49:33", "07:51", "215D_EL PATO", 62, "LP1912")
    ,@(21, "06:49:33", "08:23", "215B_EL PATO", 94, "LP1912")
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")

$ws.Range("A2").Value = "Última actualización: 06:49:33"
$ws.Range("A3").Value = "Total filas: 12"

$rows = @(
    ,@(14, "06:49:33", "07:00", "215B_LP-P MOR-1 Y 57", 11, "L6173")
    ,@(15, "06:49:33", "07:35", "215A_LA PLATA", 46, "L6173")
    ,@(16, "06:49:33", "08:07", "215C_LA PLATA", 78, "L6203")
    ,@(17, "06:49:33", "08:33", "215A_LA PLATA", 104, "L6173")
)
foreach ($row in $rows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
}
